# Reorganize scripts into 'scripts/' directory and fix path dependencies
# -> strip the leftover HTML/JS-ish scratch markup that got typed into the
#    slide bodies, and drop the now-obsolete summary slide.

$p = $ppt.ActivePresentation

# Slide 3 ("멘델의 실험 재료 완두"): the whole second paragraph was scratch
# HTML - clear it out entirely, leaving two empty paragraphs.
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$full3 = $tr3.Text
if ($full3.Length -gt 0) {
    $tr3.Characters(2, $full3.Length - 1).Delete()
}

# Slides 4-6: keep the first sentence (and the line break right after it),
# drop every extra HTML/snippet line that was appended below it.
foreach ($idx in 4, 5, 6) {
    $s = $p.Slides.Item($idx)
    $tr = $s.Shapes.Item(2).TextFrame.TextRange
    $full = $tr.Text
    $run1 = $tr.Runs(1, 1).Text
    $keepLen = $run1.Length + 1
    $deleteCount = $full.Length - $keepLen
    if ($deleteCount -gt 0) {
        $deleteStart = 1 + $keepLen + 1
        $tr.Characters($deleteStart, $deleteCount).Delete()
    }
}

# Slide 7 ("핵심 정리") is no longer needed now that content moved under
# scripts/ - remove it (also drops it from the sldIdLst automatically).
$p.Slides.Item(7).Delete()
